# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets, which hold duplicate copies of the same event list.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 6622
    6  = 2019
    7  = 1532
    9  = 1009
    10 = 423
    11 = 15
    12 = 5632
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
